$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '25.979.81'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -7.44%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.668.15'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -4.51%  '

$ws.Range('E4').Value = '  +0.71%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '217.66'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.97%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5027'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -13.60%  '

$ws.Range('E7').Value = '  +0.64%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2629'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -3.22%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06309'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -4.42%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '21.39'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -7.84%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07375'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.85%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.667.89'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -4.38%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.533'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -4.30%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5728'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -5.43%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '1.896.40'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -4.35%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.000008406'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -3.20%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '64.43'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -13.32%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '26.056.74'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -7.10%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.925'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -7.80%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.009'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.61%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.76'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -4.56%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '185.98'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -9.44%  '

$ws.Range('E23').Value = '  -7.23%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.010'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.77%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '142.75'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -4.79%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.622'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -5.11%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1162'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -6.02%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.69'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.94%  '

$ws.Range('B29').Value = 'Hedera'
$ws.Range('C29').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.05842'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -5.32%  '

$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.304'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -6.26%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.319'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -5.15%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.487'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -6.86%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.492'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -6.14%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.647'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.01%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.002'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.50%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.5965'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -6.43%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.367'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.80%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.640'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.56%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01596'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -4.73%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.080.20'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -4.15%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.931'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -4.82%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.8543'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.26%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.009'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.38%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '99.46'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.12%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.816.02'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -4.22%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00000000111'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.50%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '55.75'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -6.31%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.009'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.07%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.055'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.65%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.4314'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.40%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05178'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.73%  '
